$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new item rows above the current row 24 ("سرنجات 3 سم"),
# shifting the existing rows (and the totals/footer rows) down by two.
$ws.Rows("24:25").Insert()

# The inserted rows come in unformatted; clone the look (borders, fills,
# number formats, merges) of the data rows that used to be there and are
# now sitting two rows further down.
$ws.Range("A26:Q27").Copy()
$ws.Range("A24:Q25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 24 - new item #18
$ws.Cells.Item(24, 1).Value2 = 18
$ws.Cells.Item(24, 3).Value2 = "جل صبار برطمان"
$ws.Cells.Item(24, 8).Value2 = "5:0"
$ws.Cells.Item(24, 12).Value2 = "0"
$ws.Cells.Item(24, 14).Value2 = "25.00"
$ws.Cells.Item(24, 16).Value2 = "25.0000"
$ws.Cells.Item(24, 17).Value2 = "1:0"

# Row 25 - new item #19
$ws.Cells.Item(25, 1).Value2 = 19
$ws.Cells.Item(25, 3).Value2 = "ريكسونه رجالى"
$ws.Cells.Item(25, 8).Value2 = "5:0"
$ws.Cells.Item(25, 12).Value2 = "0"
$ws.Cells.Item(25, 14).Value2 = "27.00"
$ws.Cells.Item(25, 16).Value2 = "27.0000"
$ws.Cells.Item(25, 17).Value2 = "1:0"

# Renumber the items that used to be #18-#21 and are now two rows lower
# (#20-#23).
$ws.Cells.Item(26, 1).Value2 = 20
$ws.Cells.Item(27, 1).Value2 = 21
$ws.Cells.Item(28, 1).Value2 = 22
$ws.Cells.Item(29, 1).Value2 = 23

# The grand-total row (now row 30) grows by the price of the two new
# items (25.00 + 27.00 = 52.00).
$ws.Cells.Item(30, 16).Value2 = 891.95000000000005

# Refresh the generated-on timestamp in the footer (now row 31).
$ws.Cells.Item(31, 1).Value2 = "Wednesday, 13 August, 2025 10:49 AM"
